$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# --- Update existing rows 16 and 17 "Users" column (F) ---
$ws.Range("F16").Value = 25
$ws.Range("F17").Value = 50

# --- Fill in the newly-collected test-result columns (I:N) for rows 16-18 ---
$ws.Range("I16").Value = 93
$ws.Range("J16").Value = 118
$ws.Range("K16").Value = 1877
$ws.Range("L16").Value = 166
$ws.Range("M16").Value = 43.97
$ws.Range("N16").Value = 133.7

$ws.Range("I17").Value = 95
$ws.Range("J17").Value = 118
$ws.Range("K17").Value = 1586
$ws.Range("L17").Value = 171
$ws.Range("M17").Value = 50.6
$ws.Range("N17").Value = 260.8

$ws.Range("I18").Value = 109
$ws.Range("J18").Value = 118
$ws.Range("K18").Value = 2502
$ws.Range("L18").Value = 196
$ws.Range("M18").Value = 79.57
$ws.Range("N18").Value = 683.2

# --- Add new AWS sample rows 19-21 ---
$ws.Range("B19").Value = "AWS"
$ws.Range("C19").Value = "MathGET"
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = "NA"
$ws.Range("F19").Value = 25
$ws.Range("G19").Value = 500
$ws.Range("H19").Value = 12500
$ws.Range("I19").Value = 73
$ws.Range("J19").Value = 115
$ws.Range("K19").Value = 1167
$ws.Range("L19").Value = 127
$ws.Range("M19").Value = 31.58
$ws.Range("N19").Value = 169.7

$ws.Range("B20").Value = "AWS"
$ws.Range("C20").Value = "MathGET"
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = "NA"
$ws.Range("F20").Value = 50
$ws.Range("G20").Value = 500
$ws.Range("H20").Value = 25000
$ws.Range("I20").Value = 74
$ws.Range("J20").Value = 114
$ws.Range("K20").Value = 950
$ws.Range("L20").Value = 129
$ws.Range("M20").Value = 36.24
$ws.Range("N20").Value = 334

$ws.Range("B21").Value = "AWS"
$ws.Range("C21").Value = "MathGET"
$ws.Range("D21").Value = 1
$ws.Range("E21").Value = "NA"
$ws.Range("F21").Value = 150
$ws.Range("G21").Value = 500
$ws.Range("H21").Value = 75000
$ws.Range("I21").Value = 84
$ws.Range("J21").Value = 114
$ws.Range("K21").Value = 1998
$ws.Range("L21").Value = 141
$ws.Range("M21").Value = 66.73
$ws.Range("N21").Value = 885.9

# --- Keep the AutoFilter in sync with the new data extent (toggle off/on
#     rather than re-invoking AutoFilter() on an existing filter, which
#     would just remove it) ---
$ws.AutoFilterMode = $false
$ws.Range("B2:N21").AutoFilter() | Out-Null

# --- The hidden _FilterDatabase defined name created by AutoFilter keeps
#     pointing at the old range; repoint it at the expanded one ---
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Hoja1!_FilterDatabase") {
        $n.RefersTo = "=Hoja1!`$B`$2:`$N`$21"
    }
}

# --- Clear the stale cell selection saved in the worksheet view ---
$ws.Range("A1").Select() | Out-Null
